$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers ---
# Order of first-use matters for shared-string indices, so introduce the
# brand new header strings in this order: "Mult Speed up", "Reduction speedup",
# "Array size". The pre-existing headers for "Reduction W/O SIMD" / "Reduction
# W/ SIMD" simply move from F1/G1 to E1/F1.
$ws.Range("D1").Value = "Mult Speed up"
$ws.Range("G1").Value = "Reduction speedup"
$ws.Range("A1").Value = "Array size"
$ws.Range("E1").Value = "Reduction W/O SIMD"
$ws.Range("F1").Value = "Reduction W/ SIMD"

# --- Column D: "Mult Speed up" = C/B (replaces the old, unused B/C formula) ---
$ws.Range("D2").Formula = "=C2/B2"
$ws.Range("D3:D12").Formula = "=C3/B3"

# --- Column E: raw "array size" input data for the reduction benchmark (new) ---
$ws.Range("E2").Value2 = 128
$ws.Range("E3").Value2 = 131.24
$ws.Range("E4").Value2 = 131.87
$ws.Range("E5").Value2 = 133.15
$ws.Range("E6").Value2 = 136.88
$ws.Range("E7").Value2 = 130.5
$ws.Range("E8").Value2 = 187.86
$ws.Range("E9").Value2 = 206.55
$ws.Range("E10").Value2 = 203.85
$ws.Range("E11").Value2 = 212.12
$ws.Range("E12").Value2 = 218.42

# --- Column F: "Reduction W/ SIMD" raw values (previously held in column G) ---
$ws.Range("F2").Value2 = 881.97
$ws.Range("F3").Value2 = 976.48
$ws.Range("F4").Value2 = 996.64
$ws.Range("F5").Value2 = 990.21
$ws.Range("F6").Value2 = 982.9
$ws.Range("F7").Value2 = 969.8
$ws.Range("F8").Value2 = 998.32
$ws.Range("F9").Value2 = 993.28
$ws.Range("F10").Value2 = 979.94
$ws.Range("F11").Value2 = 976.28
$ws.Range("F12").Value2 = 1001.52

# --- Column G: new "Reduction speedup" formula column = F/E ---
$ws.Range("G2").Formula = "=F2/E2"
$ws.Range("G3:G12").Formula = "=F3/E3"

# --- Column A formulas become a shared formula block across A3:A12 ---
$ws.Range("A3:A12").Formula = "=A2+3199900"

# --- Column widths: D/E share one width, F gets its own (closest reachable
# values to the bestFit widths Excel itself would compute for this text) ---
$ws.Range("D1:E1").ColumnWidth = 17.75
$ws.Range("F1").ColumnWidth = 16.45

# --- Selection moves to E16 ---
$ws.Range("E16").Select()
